$d = $word.ActiveDocument

# The "Requisitos" bullet paragraph holds 28 runs, each "<course> - <name>  (Requisito fraco)"
# followed by a manual line break. This edit reorders those 28 lines while leaving the
# run/break structure (and the rest of the document) untouched, so we rewrite each run's
# text in place by character position rather than doing naive find/replace (which would
# collide, since the new ordering is a permutation of the existing lines).

$newLines = @(
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)",
    "LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)",
    "LOQ4095 -  Química Geral Experimental  (Requisito fraco)",
    "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)",
    "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)",
    "LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)",
    "LOB1006 -  Cálculo IV  (Requisito fraco)",
    "LOB1037 -  Àlgebra Linear  (Requisito fraco)",
    "LOB1053 -  Física III  (Requisito fraco)",
    "LOB1003 -  Cálculo I  (Requisito fraco)",
    "LOB1011 -  Eletricidade Aplicada  (Requisito fraco)",
    "LOB1012 -  Estatística  (Requisito fraco)",
    "LOB1018 -  Física I  (Requisito fraco)",
    "LOB1024 -  Mecânica  (Requisito fraco)",
    "LOB1036 -  Geometria Analítica  (Requisito fraco)",
    "LOB1038 -  Física Experimental I  (Requisito fraco)",
    "LOB1039 -  Física Experimental III  (Requisito fraco)",
    "LOB1041 -  Física Experimental II  (Requisito fraco)",
    "LOB1042 -  Física Experimental IV  (Requisito fraco)",
    "LOB1052 -  Cálculo III  (Requisito fraco)",
    "LOB1232 -  Licenciamento Ambiental  (Requisito fraco)",
    "LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)",
    "LOB1004 -  Cálculo II  (Requisito fraco)",
    "LOB1019 -  Física II  (Requisito fraco)",
    "LOB1021 -  Física IV  (Requisito fraco)",
    "LOQ4233 -  Gestão de Negócios  (Requisito fraco)"
)

# Locate the paragraph: the one styled "ListBullet" whose text contains "Requisito fraco".
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Requisito fraco*") {
        $target = $para
    }
}

$r = $target.Range
$start = $r.Start
$full = $r.Text
$segs = $full.Split([char]11)

# Compute the (start, length) of every existing run's text (segments are separated by
# the Chr(11) manual line-break character written by Word for <w:br/>).
$ranges = New-Object System.Collections.ArrayList
$pos = $start
for ($i = 0; $i -lt $segs.Count - 1; $i++) {
    $segLen = $segs[$i].Length
    [void]$ranges.Add(@($pos, $segLen))
    $pos = $pos + $segLen + 1
}

if ($ranges.Count -ne $newLines.Count) {
    throw "Expected $($newLines.Count) lines, found $($ranges.Count)"
}

# Apply replacements from the last run back to the first so earlier offsets stay valid
# even though old/new text lengths differ.
for ($i = $ranges.Count - 1; $i -ge 0; $i--) {
    $pr = $ranges[$i]
    $rng = $d.Range($pr[0], $pr[0] + $pr[1])
    $rng.Text = $newLines[$i]
}

Write-Output "Done"
